$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rx/tx header text (same cells, new wording) and add a new
# "bits" header column between rx and tx.
$ws.Range("A1").Value = "rx(что передает контроллер каналов)"
$ws.Range("C1").Value = "биты"
$ws.Range("E1").Value = "tx(что контроллер принимает)"

# Row 1 needs to grow to fit the longer wrapped header text.
$ws.Rows.Item(1).RowHeight = 45

# Widen the columns that now hold the longer header text.
$ws.Columns.Item(1).ColumnWidth = 31.83
$ws.Columns.Item(5).ColumnWidth = 28.33

# Leave the selection on C2, matching the saved view state.
$ws.Range("C2").Select()
